# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The "Estado de Cuenta" detail table (rows 16-30, columns C:G on Hoja1) is
# rebuilt: new worker records are inserted, existing worker blocks are
# resequenced, and several "Valor Mora" / "Salario Basico" amounts are
# updated to match the new source data export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row layout: TipoDoc(B) | NumDoc(C) | Nombre(D) | PeriodoMora(E) | ValorMora(F) | SalarioBasico(G)
$rows = @(
    @{ Row = 16; TipoDoc = "CC"; NumDoc = "1048938456"; Nombre = "YEVIS MARTELO MARTELO";            Periodo = "1803"; Mora = 9375;  Salario = 781242  },
    @{ Row = 17; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2301"; Mora = 12373; Salario = 1160000 },
    @{ Row = 18; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2302"; Mora = 46400; Salario = 1160000 },
    @{ Row = 19; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2303"; Mora = 46400; Salario = 1160000 },
    @{ Row = 20; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2304"; Mora = 46400; Salario = 1160000 },
    @{ Row = 21; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2305"; Mora = 46400; Salario = 1160000 },
    @{ Row = 22; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2306"; Mora = 46400; Salario = 1160000 },
    @{ Row = 23; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2307"; Mora = 46400; Salario = 1160000 },
    @{ Row = 24; TipoDoc = "CC"; NumDoc = "1044934831"; Nombre = "HERNANDO JOSE ZAMBRANO ALCALA";    Periodo = "2307"; Mora = 41765; Salario = 1423806 },
    @{ Row = 25; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2308"; Mora = 46400; Salario = 1160000 },
    @{ Row = 26; TipoDoc = "CC"; NumDoc = "1044934831"; Nombre = "HERNANDO JOSE ZAMBRANO ALCALA";    Periodo = "2308"; Mora = 56952; Salario = 1423806 },
    @{ Row = 27; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2309"; Mora = 46400; Salario = 1160000 },
    @{ Row = 28; TipoDoc = "CC"; NumDoc = "33332879";   Nombre = "LESVIS DIAZ OSPINO";                Periodo = "2309"; Mora = 4404;  Salario = 2025541 },
    @{ Row = 29; TipoDoc = "CC"; NumDoc = "1044934831"; Nombre = "HERNANDO JOSE ZAMBRANO ALCALA";    Periodo = "2309"; Mora = 1898;  Salario = 1423806 },
    @{ Row = 30; TipoDoc = "CC"; NumDoc = "1143360288"; Nombre = "JADER ENRIQUE MANJARRES OLIVERA";  Periodo = "2310"; Mora = 30934; Salario = 1160000 }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("B$n").Value = $r.TipoDoc
    $ws.Range("C$n").Value = $r.NumDoc
    $ws.Range("D$n").Value = $r.Nombre
    $ws.Range("E$n").Value = $r.Periodo
    $ws.Range("F$n").Value = $r.Mora
    $ws.Range("G$n").Value = $r.Salario
}
